# Revision History - Figma
# Aggiunte di AT e RS, modifica errori di SDS

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# reference cells that already carry the "plain text, centered" style (s=1)
# and the "centered short-date" style (s=2) used throughout the table, so
# newly written cells can pick up the same look without minting new styles
$textStyleRef = $ws.Range("E3")
$dateStyleRef = $ws.Range("B3")

# --- Row 4: Use Case Diagram Ospite (SDS) -> corrected version/description ---
$ws.Range("C4").Value = "'0.2"
$textStyleRef.Copy()
$ws.Range("C4").PasteSpecial($xlPasteFormats)
$ws.Range("D4").Value = "Aggiunta Use Case Diagram Ospite "
$ws.Range("E4").Value = "SDS"

# --- Row 5: Use Case Diagram Utente (SDS) -> corrected version/description ---
$ws.Range("C5").Value = "'0.2"
$textStyleRef.Copy()
$ws.Range("C5").PasteSpecial($xlPasteFormats)
$ws.Range("D5").Value = "Aggiunta Use Case Diagram Utente "
$ws.Range("E5").Value = "SDS"

# --- Row 6: new entry - Use Case Diagram Contadino (RS) ---
$dateStyleRef.Copy()
$ws.Range("B6").PasteSpecial($xlPasteFormats)
$ws.Range("B6").Value = "10/30/2022"

$ws.Range("C6").Value = "'0.2"
$textStyleRef.Copy()
$ws.Range("C6").PasteSpecial($xlPasteFormats)

$textStyleRef.Copy()
$ws.Range("D6").PasteSpecial($xlPasteFormats)
$ws.Range("D6").Value = "Aggiunta Use Case Diagram Contadino"

$ws.Range("E6").Value = "RS"

# --- Row 7: new entry - Use Case Diagram Catalogo (AT) ---
$dateStyleRef.Copy()
$ws.Range("B7").PasteSpecial($xlPasteFormats)
$ws.Range("B7").Value = "10/30/2022"

$ws.Range("C7").Value = "'0.2"
$textStyleRef.Copy()
$ws.Range("C7").PasteSpecial($xlPasteFormats)

$textStyleRef.Copy()
$ws.Range("D7").PasteSpecial($xlPasteFormats)
$ws.Range("D7").Value = "Aggiunta Use Case Diagram Catalogo"

$ws.Range("E7").Value = "AT"

# --- Rows 8-12: clear out the now-unused B (date) and D (description) cells ---
$ws.Range("B8").Clear()
$ws.Range("D8").Clear()
$ws.Range("B9").Clear()
$ws.Range("D9").Clear()
$ws.Range("B10").Clear()
$ws.Range("D10").Clear()
$ws.Range("B11").Clear()
$ws.Range("D11").Clear()
$ws.Range("B12").Clear()
$ws.Range("D12").Clear()

# --- Window / selection state ---
$ws.Activate()
$ws.Range("E7").Select()
$excel.ActiveWindow.Zoom = 135
